# Insert a new weekly price record at row 142, pushing the existing
# rows 142-222 down to 143-223 (dimension grows from A1:R222 to A1:R223).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(142).Insert()

$ws.Cells.Item(142, 1).Value = 4
$ws.Cells.Item(142, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(142, 3).Value = "Los Lagos"
$ws.Cells.Item(142, 4).Value = 44719
$ws.Cells.Item(142, 5).Value = 10
$ws.Cells.Item(142, 6).Value = 100112039
$ws.Cells.Item(142, 7).Value = "Ciboulette"
$ws.Cells.Item(142, 8).Value = "Sin especificar"
$ws.Cells.Item(142, 9).Value = "Primera"
$ws.Cells.Item(142, 10).Value = 240
$ws.Cells.Item(142, 11).Value = 2500
$ws.Cells.Item(142, 12).Value = 2500
$ws.Cells.Item(142, 13).Value = 2500
$ws.Cells.Item(142, 14).Value = "$/docena de atados"
$ws.Cells.Item(142, 15).Value = "Región Metropolitana"
$ws.Cells.Item(142, 16).Value = 833
$ws.Cells.Item(142, 17).Value = 3
$ws.Cells.Item(142, 18).Value = "Hortaliza"
